$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Octubre de 2020 a las 07:59"

# Row 27: Israel
$ws.Range("B27").Value = 303846
$ws.Range("C27").Value = 737
$ws.Range("D27").Value = 272015
$ws.Range("E27").Value = 29622

# Row 28: Ucrania
$ws.Range("B28").Value = 303638
$ws.Range("C28").Value = 4766
$ws.Range("D28").Value = 126489
$ws.Range("E28").Value = 171476
$ws.Range("G28").Value = 66
$ws.Range("H28").Value = 5673

# Row 62: Uzbekistan
$ws.Range("B62").Value = 63430
$ws.Range("C62").Value = 306
$ws.Range("D62").Value = 60401
$ws.Range("E62").Value = 2500
$ws.Range("G62").Value = 4
$ws.Range("H62").Value = 529

# Row 68: Kirguistan
$ws.Range("B68").Value = 52526
$ws.Range("C68").Value = 482
$ws.Range("D68").Value = 45863
$ws.Range("E68").Value = 5552

# Row 203: Santa Lucia
$ws.Range("B203").Value = 36
$ws.Range("E203").Value = 9
